# Re-order the observation rows (2-9) on the active sheet so that each
# row's species/observation data is moved to a different row, following
# the permutation cycle:
#   2 <- 3, 3 <- 9, 4 <- 5, 5 <- 6, 6 <- 2, 7 <- 8, 8 <- 7, 9 <- 4
# (new row N receives the old contents of row Source(N)).
#
# Only the cells whose value actually changes are written, so every
# untouched cell (including the already-blank "Antal"/I column cells)
# keeps its original representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value  = 111739309
$ws.Cells.Item(2,2).Value  = 78536
$ws.Cells.Item(2,4).Value  = 'LC'
$ws.Cells.Item(2,5).Value  = 229497
$ws.Cells.Item(2,6).Value  = 'Korallblylav'
$ws.Cells.Item(2,7).Value  = 'Parmeliella triptophylla'
$ws.Cells.Item(2,8).Value  = '(Ach.) Müll.Arg.'
$ws.Cells.Item(2,17).Value = 574011.1276117128
$ws.Cells.Item(2,18).Value = 7172434.078971106

$ws.Cells.Item(3,1).Value  = 111739315
$ws.Cells.Item(3,2).Value  = 78605
$ws.Cells.Item(3,5).Value  = 6462
$ws.Cells.Item(3,6).Value  = 'Stuplav'
$ws.Cells.Item(3,7).Value  = 'Nephroma bellum'
$ws.Cells.Item(3,8).Value  = '(Spreng.) Tuck.'
$ws.Cells.Item(3,17).Value = 573904.5013778479
$ws.Cells.Item(3,18).Value = 7172636.708955797

$ws.Cells.Item(4,1).Value  = 111739313
$ws.Cells.Item(4,2).Value  = 73701
$ws.Cells.Item(4,5).Value  = 1467
$ws.Cells.Item(4,6).Value  = 'Rödbrun blekspik'
$ws.Cells.Item(4,7).Value  = 'Sclerophora coniophaea'
$ws.Cells.Item(4,8).Value  = '(Norman) J.Mattsson & Middelb.'
$ws.Cells.Item(4,9).Value  = ""
$ws.Cells.Item(4,13).Value = ""
$ws.Cells.Item(4,17).Value = 574025.0565134182
$ws.Cells.Item(4,18).Value = 7172443.417908707

$ws.Cells.Item(5,1).Value  = 111739306
$ws.Cells.Item(5,2).Value  = 56398
$ws.Cells.Item(5,5).Value  = 100109
$ws.Cells.Item(5,6).Value  = 'Tretåig hackspett'
$ws.Cells.Item(5,7).Value  = 'Picoides tridactylus'
$ws.Cells.Item(5,8).Value  = '(Linnaeus, 1758)'
$ws.Cells.Item(5,13).Value = 'äldre spår'
$ws.Cells.Item(5,17).Value = 573906.0397215446
$ws.Cells.Item(5,18).Value = 7172521.061635921

$ws.Cells.Item(6,1).Value  = 111739316
$ws.Cells.Item(6,2).Value  = 78578
$ws.Cells.Item(6,5).Value  = 6458
$ws.Cells.Item(6,6).Value  = 'Lunglav'
$ws.Cells.Item(6,7).Value  = 'Lobaria pulmonaria'
$ws.Cells.Item(6,8).Value  = '(L.) Hoffm.'
$ws.Cells.Item(6,13).Value = ""
$ws.Cells.Item(6,17).Value = 573904.5013778479
$ws.Cells.Item(6,18).Value = 7172636.708955797

$ws.Cells.Item(7,1).Value  = 111739317
$ws.Cells.Item(7,2).Value  = 78579
$ws.Cells.Item(7,5).Value  = 2081
$ws.Cells.Item(7,6).Value  = 'Skrovellav'
$ws.Cells.Item(7,7).Value  = 'Lobaria scrobiculata'
$ws.Cells.Item(7,8).Value  = '(Scop.) DC.'
$ws.Cells.Item(7,17).Value = 573911.5177193542
$ws.Cells.Item(7,18).Value = 7172648.020174325

$ws.Cells.Item(8,1).Value  = 111739311
$ws.Cells.Item(8,2).Value  = 77515
$ws.Cells.Item(8,5).Value  = 6425
$ws.Cells.Item(8,6).Value  = 'Garnlav'
$ws.Cells.Item(8,7).Value  = 'Alectoria sarmentosa'
$ws.Cells.Item(8,8).Value  = '(Ach.) Ach.'
$ws.Cells.Item(8,17).Value = 574011.8892867711
$ws.Cells.Item(8,18).Value = 7172473.089384713

$ws.Cells.Item(9,1).Value  = 111739307
$ws.Cells.Item(9,2).Value  = 56543
$ws.Cells.Item(9,4).Value  = 'NT'
$ws.Cells.Item(9,5).Value  = 103021
$ws.Cells.Item(9,6).Value  = 'Talltita'
$ws.Cells.Item(9,7).Value  = 'Poecile montanus'
$ws.Cells.Item(9,8).Value  = '(Conrad von Baldenstein, 1827)'
$ws.Cells.Item(9,9).Value  = '3'
$ws.Cells.Item(9,13).Value = 'födosökande'
$ws.Cells.Item(9,17).Value = 573960.5743707293
$ws.Cells.Item(9,18).Value = 7172501.399265604
